$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 94.22221999999999
$ws.Range("I5").Value = 49.833332
$ws.Range("J5").Value = 183
$ws.Range("K5").Value = 49.833332
$ws.Range("L5").Value = 183
$ws.Range("M5").Value = 65.166668
$ws.Range("N5").Value = -413
$ws.Range("H18").Value = 1070.4117
$ws.Range("I18").Value = 762.3125
$ws.Range("J18").Value = 6000
$ws.Range("K18").Value = 762.3125
$ws.Range("L18").Value = 6000
$ws.Range("M18").Value = -478.3125
$ws.Range("N18").Value = -6568
$ws.Range("H64").Value = 37269.45
$ws.Range("J64").Value = 2986.3333
$ws.Range("L64").Value = 2986.3333
$ws.Range("N64").Value = -3482.3333
$ws.Range("H67").Value = 37269.45
$ws.Range("J67").Value = 2986.3333
$ws.Range("L67").Value = 2986.3333
$ws.Range("N67").Value = -4702.3333
$ws.Range("H74").Value = 3566.9167
$ws.Range("I74").Value = 3200.3333
$ws.Range("K74").Value = 3200.3333
$ws.Range("M74").Value = -2264.3333
$ws.Range("H77").Value = 3566.9167
$ws.Range("I77").Value = 3200.3333
$ws.Range("K77").Value = 16001.6665
$ws.Range("M77").Value = -11321.6665
$ws.Range("H129").Value = 3534.8538
$ws.Range("I129").Value = 8892.666999999999
$ws.Range("J129").Value = 1317.8276
$ws.Range("K129").Value = 26678.001
$ws.Range("L129").Value = 3953.4828
$ws.Range("M129").Value = -21678.001
$ws.Range("N129").Value = -13953.4828
$ws.Range("H132").Value = 4469997
$ws.Range("I132").Value = 5214230
$ws.Range("J132").Value = 4597.375
$ws.Range("K132").Value = 15642690
$ws.Range("L132").Value = 13792.125
$ws.Range("M132").Value = -15640160
$ws.Range("N132").Value = -18852.125
$ws.Range("H138").Value = 2518.0266
$ws.Range("I138").Value = 1277.4849
$ws.Range("J138").Value = 3492.738
$ws.Range("K138").Value = 3832.4547
$ws.Range("L138").Value = 10478.214
$ws.Range("M138").Value = 1307.5453
$ws.Range("N138").Value = -20758.214

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 38186.203
$ws.Range("I32").Value = 10792.2295
$ws.Range("J32").Value = 257338
$ws.Range("K32").Value = 10792.2295
$ws.Range("L32").Value = 257338
$ws.Range("M32").Value = -10505.2295
$ws.Range("N32").Value = -257912
$ws.Range("H74").Value = 829.6667
$ws.Range("I74").Value = 718
$ws.Range("J74").Value = 1164.6666
$ws.Range("K74").Value = 718
$ws.Range("L74").Value = 1164.6666
$ws.Range("M74").Value = 156
$ws.Range("N74").Value = -2912.6666
$ws.Range("H77").Value = 829.6667
$ws.Range("I77").Value = 718
$ws.Range("J77").Value = 1164.6666
$ws.Range("K77").Value = 3590
$ws.Range("L77").Value = 5823.333000000001
$ws.Range("M77").Value = 778
$ws.Range("N77").Value = -14559.333
$ws.Range("H110").Value = 125263220
$ws.Range("I110").Value = 125263220
$ws.Range("K110").Value = 125263220
$ws.Range("M110").Value = -125261175
$ws.Range("H122").Value = 2111.9312
$ws.Range("I122").Value = 1740.3636
$ws.Range("K122").Value = 5221.0908
$ws.Range("M122").Value = -2771.0908
$ws.Range("H132").Value = 13580.551
$ws.Range("I132").Value = 16246.538
$ws.Range("J132").Value = 3183.2
$ws.Range("K132").Value = 48739.614
$ws.Range("L132").Value = 9549.599999999999
$ws.Range("M132").Value = -46209.614
$ws.Range("N132").Value = -14609.6

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 18541
$ws.Range("J21").Value = 18541
$ws.Range("L21").Value = 18541
$ws.Range("N21").Value = -19013
$ws.Range("H86").Value = 66881.3
$ws.Range("I86").Value = 86523.234
$ws.Range("J86").Value = 3045
$ws.Range("K86").Value = 86523.234
$ws.Range("L86").Value = 3045
$ws.Range("M86").Value = -85400.234
$ws.Range("N86").Value = -5291
$ws.Range("H89").Value = 66881.3
$ws.Range("I89").Value = 86523.234
$ws.Range("J89").Value = 3045
$ws.Range("K89").Value = 432616.17
$ws.Range("L89").Value = 15225
$ws.Range("M89").Value = -427000.17
$ws.Range("N89").Value = -26457

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 1664.2858
$ws.Range("J8").Value = 2010
$ws.Range("L8").Value = 2010
$ws.Range("N8").Value = -2290
$ws.Range("H17").Value = 3835.6
$ws.Range("I17").Value = 294.5
$ws.Range("K17").Value = 294.5
$ws.Range("M17").Value = -120.5
$ws.Range("H22").Value = 320.07693
$ws.Range("I22").Value = 236.1
$ws.Range("K22").Value = 236.1
$ws.Range("M22").Value = 113.9
$ws.Range("H25").Value = 19900
$ws.Range("J25").Value = 19900
$ws.Range("L25").Value = 19900
$ws.Range("N25").Value = -20248
$ws.Range("H41").Value = 9049.1
$ws.Range("I41").Value = 4310.2
$ws.Range("J41").Value = 13788
$ws.Range("K41").Value = 4310.2
$ws.Range("L41").Value = 13788
$ws.Range("M41").Value = -3882.2
$ws.Range("N41").Value = -14644
$ws.Range("H50").Value = 13990
$ws.Range("J50").Value = 13990
$ws.Range("L50").Value = 13990
$ws.Range("N50").Value = -15240
$ws.Range("H58").Value = 1436.3158
$ws.Range("I58").Value = 1188.6923
$ws.Range("J58").Value = 1972.8334
$ws.Range("K58").Value = 1188.6923
$ws.Range("L58").Value = 1972.8334
$ws.Range("M58").Value = -985.6922999999999
$ws.Range("N58").Value = -2378.8334
$ws.Range("H60").Value = 13021.818
$ws.Range("I60").Value = 10424.25
$ws.Range("J60").Value = 14506.143
$ws.Range("K60").Value = 10424.25
$ws.Range("L60").Value = 14506.143
$ws.Range("M60").Value = -9913.25
$ws.Range("N60").Value = -15528.143
$ws.Range("H81").Value = 37098
$ws.Range("J81").Value = 37098
$ws.Range("L81").Value = 37098
$ws.Range("N81").Value = -39094
$ws.Range("H84").Value = 37098
$ws.Range("J84").Value = 37098
$ws.Range("L84").Value = 111294
$ws.Range("N84").Value = -121278
$ws.Range("H136").Value = 1436.3158
$ws.Range("I136").Value = 1188.6923
$ws.Range("J136").Value = 1972.8334
$ws.Range("K136").Value = 3566.0769
$ws.Range("L136").Value = 5918.5002
$ws.Range("M136").Value = -1016.0769
$ws.Range("N136").Value = -11018.5002

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 613.1667
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 9000
$ws.Range("N122").Value = -13900

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 143173140
$ws.Range("I80").Value = 167034990
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 167034990
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -167033992
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 143173140
$ws.Range("I83").Value = 167034990
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 835174950
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -835169958
$ws.Range("N83").Value = -19984
$ws.Range("H97").Value = 40002256
$ws.Range("I97").Value = 52633916
$ws.Range("J97").Value = 2003.6666
$ws.Range("K97").Value = 52633916
$ws.Range("L97").Value = 2003.6666
$ws.Range("M97").Value = -52633420
$ws.Range("N97").Value = -2995.6666
$ws.Range("H102").Value = 1942
$ws.Range("I102").Value = 1654.3636
$ws.Range("J102").Value = 2185.3845
$ws.Range("K102").Value = 1654.3636
$ws.Range("L102").Value = 2185.3845
$ws.Range("M102").Value = -32.36359999999991
$ws.Range("N102").Value = -5429.3845
$ws.Range("H122").Value = 1859.28
$ws.Range("I122").Value = 1475.7273
$ws.Range("J122").Value = 4672
$ws.Range("K122").Value = 4427.1819
$ws.Range("L122").Value = 14016
$ws.Range("M122").Value = -1977.1819
$ws.Range("N122").Value = -18916
$ws.Range("H132").Value = 3558
$ws.Range("I132").Value = 2687.1765
$ws.Range("J132").Value = 5672.857
$ws.Range("K132").Value = 8061.529500000001
$ws.Range("L132").Value = 17018.571
$ws.Range("M132").Value = -5531.529500000001
$ws.Range("N132").Value = -22078.571

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 608.875
$ws.Range("I55").Value = 917.2381
$ws.Range("J55").Value = 369.03705
$ws.Range("K55").Value = 917.2381
$ws.Range("L55").Value = 369.03705
$ws.Range("M55").Value = -744.2381
$ws.Range("N55").Value = -715.03705
$ws.Range("H82").Value = 2231.5
$ws.Range("I82").Value = 1775.5
$ws.Range("K82").Value = 1775.5
$ws.Range("M82").Value = -1414.5
$ws.Range("H85").Value = 2231.5
$ws.Range("I85").Value = 1775.5
$ws.Range("K85").Value = 1775.5
$ws.Range("M85").Value = -527.5
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990
$ws.Range("H100").Value = 2239.5715
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 2335.4
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 2335.4
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -3417.4
$ws.Range("H132").Value = 2957.3667
$ws.Range("I132").Value = 3167.4443
$ws.Range("K132").Value = 9502.332900000001
$ws.Range("M132").Value = -6972.332900000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 40000
$ws.Range("J97").Value = 40000
$ws.Range("L97").Value = 40000
$ws.Range("N97").Value = -41982
$ws.Range("H126").Value = 1462.4762
$ws.Range("I126").Value = 1547.3125
$ws.Range("K126").Value = 4641.9375
$ws.Range("M126").Value = -2171.9375
$ws.Range("H132").Value = 5134.073
$ws.Range("I132").Value = 3182.64
$ws.Range("J132").Value = 8183.1875
$ws.Range("K132").Value = 9547.92
$ws.Range("L132").Value = 24549.5625
$ws.Range("M132").Value = -7017.92
$ws.Range("N132").Value = -29609.5625
$ws.Range("H136").Value = 20322.186
$ws.Range("I136").Value = 48252.57
$ws.Range("J136").Value = 4886.9736
$ws.Range("K136").Value = 144757.71
$ws.Range("L136").Value = 14660.9208
$ws.Range("M136").Value = -142207.71
$ws.Range("N136").Value = -19760.9208
